# Applies the text replacements described by the diff:
# date header + 25 two-digit multiplication problems in the table.
$d = $word.ActiveDocument
$range = $d.Content

$found = $range.Find.Execute("2023-12-31 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-01 Monday", 2)
if (-not $found) { throw "Could not find text: 2023-12-31 Sunday" }
$found = $range.Find.Execute("33×95=3135", $true, $false, $false, $false, $false, $true, 1, $false, "86×52=4472", 2)
if (-not $found) { throw "Could not find text: 33×95=3135" }
$found = $range.Find.Execute("98×18=1764", $true, $false, $false, $false, $false, $true, 1, $false, "34×72=2448", 2)
if (-not $found) { throw "Could not find text: 98×18=1764" }
$found = $range.Find.Execute("20×27=540", $true, $false, $false, $false, $false, $true, 1, $false, "73×62=4526", 2)
if (-not $found) { throw "Could not find text: 20×27=540" }
$found = $range.Find.Execute("38×39=1482", $true, $false, $false, $false, $false, $true, 1, $false, "70×61=4270", 2)
if (-not $found) { throw "Could not find text: 38×39=1482" }
$found = $range.Find.Execute("14×93=1302", $true, $false, $false, $false, $false, $true, 1, $false, "88×24=2112", 2)
if (-not $found) { throw "Could not find text: 14×93=1302" }
$found = $range.Find.Execute("52×74=3848", $true, $false, $false, $false, $false, $true, 1, $false, "96×94=9024", 2)
if (-not $found) { throw "Could not find text: 52×74=3848" }
$found = $range.Find.Execute("48×50=2400", $true, $false, $false, $false, $false, $true, 1, $false, "76×90=6840", 2)
if (-not $found) { throw "Could not find text: 48×50=2400" }
$found = $range.Find.Execute("13×95=1235", $true, $false, $false, $false, $false, $true, 1, $false, "68×84=5712", 2)
if (-not $found) { throw "Could not find text: 13×95=1235" }
$found = $range.Find.Execute("40×96=3840", $true, $false, $false, $false, $false, $true, 1, $false, "74×27=1998", 2)
if (-not $found) { throw "Could not find text: 40×96=3840" }
$found = $range.Find.Execute("25×23=575", $true, $false, $false, $false, $false, $true, 1, $false, "13×55=715", 2)
if (-not $found) { throw "Could not find text: 25×23=575" }
$found = $range.Find.Execute("34×68=2312", $true, $false, $false, $false, $false, $true, 1, $false, "25×36=900", 2)
if (-not $found) { throw "Could not find text: 34×68=2312" }
$found = $range.Find.Execute("53×18=954", $true, $false, $false, $false, $false, $true, 1, $false, "83×54=4482", 2)
if (-not $found) { throw "Could not find text: 53×18=954" }
$found = $range.Find.Execute("37×23=851", $true, $false, $false, $false, $false, $true, 1, $false, "57×74=4218", 2)
if (-not $found) { throw "Could not find text: 37×23=851" }
$found = $range.Find.Execute("77×80=6160", $true, $false, $false, $false, $false, $true, 1, $false, "81×85=6885", 2)
if (-not $found) { throw "Could not find text: 77×80=6160" }
$found = $range.Find.Execute("48×85=4080", $true, $false, $false, $false, $false, $true, 1, $false, "38×77=2926", 2)
if (-not $found) { throw "Could not find text: 48×85=4080" }
$found = $range.Find.Execute("29×24=696", $true, $false, $false, $false, $false, $true, 1, $false, "34×37=1258", 2)
if (-not $found) { throw "Could not find text: 29×24=696" }
$found = $range.Find.Execute("96×30=2880", $true, $false, $false, $false, $false, $true, 1, $false, "53×97=5141", 2)
if (-not $found) { throw "Could not find text: 96×30=2880" }
$found = $range.Find.Execute("32×22=704", $true, $false, $false, $false, $false, $true, 1, $false, "46×98=4508", 2)
if (-not $found) { throw "Could not find text: 32×22=704" }
$found = $range.Find.Execute("33×51=1683", $true, $false, $false, $false, $false, $true, 1, $false, "34×18=612", 2)
if (-not $found) { throw "Could not find text: 33×51=1683" }
$found = $range.Find.Execute("63×41=2583", $true, $false, $false, $false, $false, $true, 1, $false, "43×21=903", 2)
if (-not $found) { throw "Could not find text: 63×41=2583" }
$found = $range.Find.Execute("89×24=2136", $true, $false, $false, $false, $false, $true, 1, $false, "95×77=7315", 2)
if (-not $found) { throw "Could not find text: 89×24=2136" }
$found = $range.Find.Execute("35×97=3395", $true, $false, $false, $false, $false, $true, 1, $false, "65×90=5850", 2)
if (-not $found) { throw "Could not find text: 35×97=3395" }
$found = $range.Find.Execute("22×42=924", $true, $false, $false, $false, $false, $true, 1, $false, "26×60=1560", 2)
if (-not $found) { throw "Could not find text: 22×42=924" }
$found = $range.Find.Execute("76×44=3344", $true, $false, $false, $false, $false, $true, 1, $false, "53×27=1431", 2)
if (-not $found) { throw "Could not find text: 76×44=3344" }
$found = $range.Find.Execute("54×37=1998", $true, $false, $false, $false, $false, $true, 1, $false, "46×75=3450", 2)
if (-not $found) { throw "Could not find text: 54×37=1998" }
